$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '45.183.35'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +3.32%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.366.04'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +1.15%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.17%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '310.65'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.28%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '107.86'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -1.41%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.13%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.616'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -0.54%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '40.81'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -0.78%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.56%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '8.45'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -1.68%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +1.28%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.975'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -3.02%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.726.52'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.51%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '15.23'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -2.04%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.384.83'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +2.36%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '45.129.88'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +3.49%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '14.62'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +10.93%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.27'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -4.44%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -1.22%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '73.12'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -1.64%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.70%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '259.30'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -3.85%  '
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.28%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.24%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.09'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.79%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.23'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -5.37%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.35'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0965'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +8.69%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -1.57%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '37.32'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -3.82%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '168.93'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +0.60%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +5.88%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -1.07%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +3.38%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -1.08%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +2.68%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.90'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +2.04%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0352'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -3.89%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.75'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +1.00%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '99.42'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -5.41%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -3.53%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '69.24'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -3.59%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '12.93'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -3.11%  '
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.18%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.848.63'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +10.83%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '81.59'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +5.56%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '5.58'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +4.56%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '111.74'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -2.38%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '9.14'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +1.95%  '
